$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update betting odds values per the 2025-05-22 FlashScore refresh.
$ws.Cells.Item(10, 7).Value = 1.65  # G10: 1.7 -> 1.65
$ws.Cells.Item(10, 8).Value = 4.2  # H10: 4 -> 4.2
$ws.Cells.Item(10, 20).Value = 8.5  # T10: 8 -> 8.5
$ws.Cells.Item(10, 21).Value = 9  # U10: 8.5 -> 9
$ws.Cells.Item(10, 25).Value = 21  # Y10: 23 -> 21
$ws.Cells.Item(10, 30).Value = 151  # AD10: 201 -> 151
$ws.Cells.Item(10, 32).Value = 26  # AF10: 23 -> 26
$ws.Cells.Item(12, 9).Value = 6  # I12: 6.25 -> 6
$ws.Cells.Item(12, 23).Value = 11  # W12: 10 -> 11
$ws.Cells.Item(13, 8).Value = 3.6  # H13: 3.5 -> 3.6
$ws.Cells.Item(13, 9).Value = 4  # I13: 3.9 -> 4
$ws.Cells.Item(13, 10).Value = 1.06  # J13: 1.05 -> 1.06
$ws.Cells.Item(13, 11).Value = 10  # K13: 11 -> 10
$ws.Cells.Item(13, 26).Value = 10  # Z13: 11 -> 10
$ws.Cells.Item(13, 31).Value = 11  # AE13: 12 -> 11
$ws.Cells.Item(14, 10).Value = 1.05  # J14: 1.06 -> 1.05
$ws.Cells.Item(14, 11).Value = 11  # K14: 10 -> 11
$ws.Cells.Item(15, 10).Value = 1.04  # J15: 1.02 -> 1.04
$ws.Cells.Item(15, 11).Value = 9  # K15: 12 -> 9
$ws.Cells.Item(17, 10).Value = 1.02  # J17: 19 -> 1.02
$ws.Cells.Item(17, 11).Value = 11  # K17: 1.03 -> 11
$ws.Cells.Item(19, 7).Value = 1.7  # G19: 1.75 -> 1.7
$ws.Cells.Item(19, 8).Value = 3.9  # H19: 3.8 -> 3.9
$ws.Cells.Item(19, 9).Value = 4.75  # I19: 4.33 -> 4.75
$ws.Cells.Item(19, 24).Value = 13  # X19: 15 -> 13
$ws.Cells.Item(19, 27).Value = 7.5  # AA19: 7 -> 7.5
$ws.Cells.Item(19, 35).Value = 41  # AI19: 34 -> 41
$ws.Cells.Item(20, 10).Value = 1.06  # J20: 1.05 -> 1.06
$ws.Cells.Item(20, 11).Value = 10  # K20: 11 -> 10
$ws.Cells.Item(20, 14).Value = 2.03  # N20: 2 -> 2.03
$ws.Cells.Item(20, 15).Value = 1.78  # O20: 1.8 -> 1.78
$ws.Cells.Item(21, 9).Value = 3.4  # I21: 3.5 -> 3.4
$ws.Cells.Item(21, 10).Value = 1.05  # J21: 1.06 -> 1.05
$ws.Cells.Item(21, 11).Value = 11  # K21: 10 -> 11
$ws.Cells.Item(21, 12).Value = 1.3  # L21: 1.33 -> 1.3
$ws.Cells.Item(21, 13).Value = 3.4  # M21: 3.25 -> 3.4
$ws.Cells.Item(21, 14).Value = 2  # N21: 2.05 -> 2
$ws.Cells.Item(21, 15).Value = 1.8  # O21: 1.75 -> 1.8
$ws.Cells.Item(21, 18).Value = 1.8  # R21: 1.91 -> 1.8
$ws.Cells.Item(21, 19).Value = 1.95  # S21: 1.91 -> 1.95
$ws.Cells.Item(21, 26).Value = 10  # Z21: 9.5 -> 10
$ws.Cells.Item(21, 31).Value = 10  # AE21: 9.5 -> 10
$ws.Cells.Item(21, 36).Value = 34  # AJ21: 41 -> 34
$ws.Cells.Item(23, 12).Value = 1.3  # L23: 1.33 -> 1.3
$ws.Cells.Item(23, 13).Value = 3.4  # M23: 3.25 -> 3.4
$ws.Cells.Item(23, 14).Value = 2.05  # N23: 2.08 -> 2.05
$ws.Cells.Item(23, 15).Value = 1.75  # O23: 1.73 -> 1.75
$ws.Cells.Item(27, 7).Value = 2.05  # G27: 1.95 -> 2.05
$ws.Cells.Item(27, 8).Value = 3.7  # H27: 3.75 -> 3.7
$ws.Cells.Item(27, 9).Value = 3.2  # I27: 3.4 -> 3.2
$ws.Cells.Item(27, 18).Value = 1.5  # R27: 1.53 -> 1.5
$ws.Cells.Item(27, 19).Value = 2.5  # S27: 2.38 -> 2.5
$ws.Cells.Item(27, 20).Value = 11  # T27: 10 -> 11
$ws.Cells.Item(27, 21).Value = 12  # U27: 11 -> 12
$ws.Cells.Item(27, 26).Value = 17  # Z27: 15 -> 17
$ws.Cells.Item(27, 30).Value = 101  # AD27: 126 -> 101
$ws.Cells.Item(27, 31).Value = 13  # AE27: 15 -> 13
$ws.Cells.Item(27, 32).Value = 19  # AF27: 21 -> 19
$ws.Cells.Item(27, 34).Value = 34  # AH27: 41 -> 34
$ws.Cells.Item(29, 7).Value = 2.55  # G29: 2.57 -> 2.55
$ws.Cells.Item(29, 9).Value = 3.05  # I29: 3 -> 3.05
$ws.Cells.Item(29, 12).Value = 1.5  # L29: 1.52 -> 1.5
$ws.Cells.Item(29, 13).Value = 2.27  # M29: 2.22 -> 2.27
$ws.Cells.Item(29, 14).Value = 2.42  # N29: 2.47 -> 2.42
$ws.Cells.Item(29, 15).Value = 1.44  # O29: 1.42 -> 1.44
$ws.Cells.Item(29, 16).Value = 1.53  # P29: 1.55 -> 1.53
$ws.Cells.Item(29, 17).Value = 2.18  # Q29: 2.15 -> 2.18
$ws.Cells.Item(29, 18).Value = 1.98  # R29: 2.02 -> 1.98
$ws.Cells.Item(29, 19).Value = 1.65  # S29: 1.62 -> 1.65
$ws.Cells.Item(29, 21).Value = 11.25  # U29: 11.5 -> 11.25
$ws.Cells.Item(29, 22).Value = 10  # V29: 10.25 -> 10
$ws.Cells.Item(29, 23).Value = 29  # W29: 30 -> 29
$ws.Cells.Item(29, 24).Value = 26  # X29: 27 -> 26
$ws.Cells.Item(29, 26).Value = 6  # Z29: 5.9 -> 6
$ws.Cells.Item(29, 28).Value = 16.5  # AB29: 17 -> 16.5
$ws.Cells.Item(29, 31).Value = 7.1  # AE29: 6.8 -> 7.1
$ws.Cells.Item(29, 32).Value = 14.5  # AF29: 14 -> 14.5
$ws.Cells.Item(29, 33).Value = 11.25  # AG29: 11.5 -> 11.25
$ws.Cells.Item(29, 35).Value = 32  # AI29: 35 -> 32
$ws.Cells.Item(29, 36).Value = 50  # AJ29: 55 -> 50
$ws.Cells.Item(30, 7).Value = 2.12  # G30: 2.15 -> 2.12
$ws.Cells.Item(30, 8).Value = 3.1  # H30: 3 -> 3.1
$ws.Cells.Item(30, 9).Value = 3.4  # I30: 3.45 -> 3.4
$ws.Cells.Item(30, 12).Value = 1.35  # L30: 1.36 -> 1.35
$ws.Cells.Item(30, 13).Value = 2.7  # M30: 2.65 -> 2.7
$ws.Cells.Item(30, 14).Value = 2.02  # N30: 2.05 -> 2.02
$ws.Cells.Item(30, 15).Value = 1.62  # O30: 1.6 -> 1.62
$ws.Cells.Item(30, 19).Value = 1.82  # S30: 1.8 -> 1.82
$ws.Cells.Item(30, 20).Value = 6.9  # T30: 6.5 -> 6.9
$ws.Cells.Item(30, 21).Value = 10  # U30: 9.75 -> 10
$ws.Cells.Item(30, 23).Value = 20  # W30: 21 -> 20
$ws.Cells.Item(30, 24).Value = 18  # X30: 19 -> 18
$ws.Cells.Item(30, 25).Value = 30  # Y30: 32 -> 30
$ws.Cells.Item(30, 26).Value = 8.25  # Z30: 7.9 -> 8.25
$ws.Cells.Item(30, 27).Value = 6  # AA30: 5.9 -> 6
$ws.Cells.Item(30, 28).Value = 15  # AB30: 14.5 -> 15
$ws.Cells.Item(30, 31).Value = 9  # AE30: 9.25 -> 9
$ws.Cells.Item(30, 32).Value = 17.5  # AF30: 18 -> 17.5
$ws.Cells.Item(30, 36).Value = 45  # AJ30: 40 -> 45
